$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param(
        $ws,
        [int]$row,
        [hashtable]$values
    )
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

# ---- ALC (sheet1) ----
$ws = $wb.Worksheets.Item("ALC")

Set-RowValues $ws 64 @{
    H = 3013.0557
    I = 2761.111
    J = 3265
    K = 2761.111
    L = 3265
    M = -2513.111
    N = -3761
}

Set-RowValues $ws 67 @{
    H = 3013.0557
    I = 2761.111
    J = 3265
    K = 2761.111
    L = 3265
    M = -1903.111
    N = -4981
}

Set-RowValues $ws 76 @{
    H = 2937.6
    I = 2723.0588
    K = 2723.0588
    M = -2408.0588
}

Set-RowValues $ws 79 @{
    H = 2937.6
    I = 2723.0588
    K = 2723.0588
    M = -1631.0588
}

Set-RowValues $ws 132 @{
    H = 9299.264999999999
    I = 9946.315000000001
    K = 29838.945
    M = -27308.945
}

# ---- ARM (sheet2) ----
$ws = $wb.Worksheets.Item("ARM")

Set-RowValues $ws 63 @{
    H = 2930.6448
    I = 2235.6604
    J = 4532.1304
    K = 2235.6604
    L = 4532.1304
    M = -1549.6604
    N = -5904.1304
}

Set-RowValues $ws 66 @{
    H = 2930.6448
    I = 2235.6604
    J = 4532.1304
    K = 11178.302
    L = 22660.652
    M = -7746.302000000001
    N = -29524.652
}

Set-RowValues $ws 88 @{
    H = 1896.0769
    I = 1791.6666
    J = 1985.5714
    K = 1791.6666
    L = 1985.5714
    M = -1385.6666
    N = -2797.5714
}

Set-RowValues $ws 91 @{
    H = 1896.0769
    I = 1791.6666
    J = 1985.5714
    K = 1791.6666
    L = 1985.5714
    M = -387.6666
    N = -4793.5714
}

# ---- BSM (sheet3) ----
$ws = $wb.Worksheets.Item("BSM")

Set-RowValues $ws 86 @{
    H = 1513.4736
    I = 1496.3334
    J = 1542.8572
    K = 1496.3334
    L = 1542.8572
    M = -373.3334
    N = -3788.8572
}

Set-RowValues $ws 89 @{
    H = 1513.4736
    I = 1496.3334
    J = 1542.8572
    K = 7481.666999999999
    L = 7714.286
    M = -1865.666999999999
    N = -18946.286
}

Set-RowValues $ws 105 @{
    H = 2041.1538
    I = 1828.5555
    K = 1828.5555
    M = -81.55549999999994
}

# ---- CRP (sheet4) ----
$ws = $wb.Worksheets.Item("CRP")

Set-RowValues $ws 7 @{
    H = 2681.2632
    I = 10063.2
    J = 44.857143
    K = 10063.2
    L = 44.857143
    M = -9950.200000000001
    N = -270.857143
}

Set-RowValues $ws 31 @{
    H = 3877997
    I = 1265.5178
    K = 1265.5178
    M = -970.5178000000001
}

Set-RowValues $ws 34 @{
    H = 3877997
    I = 1265.5178
    K = 1265.5178
    M = -1063.5178
}

Set-RowValues $ws 43 @{
    H = 21551.4
    J = 21551.4
    L = 21551.4
    N = -21919.4
}

Set-RowValues $ws 62 @{
    H = 12838.5
    J = 22738
    L = 22738
    N = -23986
}

Set-RowValues $ws 65 @{
    H = 12838.5
    J = 22738
    L = 113690
    N = -119930
}

Set-RowValues $ws 94 @{
    H = 1617
    I = 1234
    K = 1234
    M = -783
}

Set-RowValues $ws 99 @{
    H = 2659
    I = 2442.8572
    J = 3163.3333
    K = 2442.8572
    L = 3163.3333
    M = -944.8571999999999
    N = -6159.3333
}

Set-RowValues $ws 101 @{
    H = 21551.4
    J = 21551.4
    L = 21551.4
    N = -28041.4
}

Set-RowValues $ws 126 @{
    H = 2659
    I = 2442.8572
    J = 3163.3333
    K = 7328.571599999999
    L = 9489.999899999999
    M = -4858.571599999999
    N = -14429.9999
}

Set-RowValues $ws 134 @{
    H = 4876.8
    I = 2834.6667
    J = 7940
    K = 8504.000100000001
    L = 23820
    M = -5969.000100000001
    N = -28890
}

# ---- GSM (sheet6) ----
$ws = $wb.Worksheets.Item("GSM")

Set-RowValues $ws 70 @{
    H = 6258.85
    I = 6360.2085
    J = 6106.8125
    K = 6360.2085
    L = 6106.8125
    M = -6090.2085
    N = -6646.8125
}

Set-RowValues $ws 73 @{
    H = 6258.85
    I = 6360.2085
    J = 6106.8125
    K = 6360.2085
    L = 6106.8125
    M = -5424.2085
    N = -7978.8125
}

Set-RowValues $ws 80 @{
    H = 2491.6365
    I = 2500.8333
    J = 2480.6
    K = 2500.8333
    L = 2480.6
    M = -1502.8333
    N = -4476.6
}

Set-RowValues $ws 83 @{
    H = 2491.6365
    I = 2500.8333
    J = 2480.6
    K = 12504.1665
    L = 12403
    M = -7512.166499999999
    N = -22387
}

# ---- LTW (sheet7) ----
$ws = $wb.Worksheets.Item("LTW")

Set-RowValues $ws 103 @{
    H = 78301
    J = 78301
    L = 78301
    N = -80645
}

# ---- WVR (sheet8) ----
$ws = $wb.Worksheets.Item("WVR")

Set-RowValues $ws 103 @{
    H = 0
    J = 0
    L = 0
}
$ws.Range("N103").ClearContents()
